$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 23 to the "Tipi di mutazione" table:
#   A23 = 22 (code, continuing the existing sequence)
#   B23 = "Dati decesso" (new description)
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "Dati decesso"

# Give the new row the same formatting as the rest of the data rows (row 22)
$ws.Range("A22:B22").Copy()
$ws.Range("A23:B23").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Match the author's final selection state
$ws.Range("B23").Select()
